$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.606.25'
$ws.Range("E2").Value = '  +3.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.256.68'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.14'
$ws.Range("E5").Value = '  +2.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '90.91'
$ws.Range("E6").Value = '  +3.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  +3.23%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  +1.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.00'
$ws.Range("E10").Value = '  +3.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.84'
$ws.Range("E11").Value = '  +2.80%  '

$ws.Range("E12").Value = '  +1.76%  '

$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("E14").Value = '  +2.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.606.51'
$ws.Range("E15").Value = '  +2.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.12'
$ws.Range("E16").Value = '  +2.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.249.96'
$ws.Range("E17").Value = '  -5.92%  '

$ws.Range("E18").Value = '  +2.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.538.94'
$ws.Range("E19").Value = '  +3.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.37'
$ws.Range("E20").Value = '  +10.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("E22").Value = '  +2.44%  '

$ws.Range("E23").Value = '  +1.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.29'
$ws.Range("E24").Value = '  +2.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  +3.89%  '

$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("E27").Value = '  +5.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.89'
$ws.Range("E28").Value = '  +3.14%  '

$ws.Range("E29").Value = '  +1.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '161.91'
$ws.Range("E31").Value = '  +1.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.24'
$ws.Range("E32").Value = '  +7.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.14'
$ws.Range("E34").Value = '  +3.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0740'
$ws.Range("E35").Value = '  +3.63%  '

$ws.Range("E36").Value = '  -1.54%  '

$ws.Range("E37").Value = '  +2.13%  '

$ws.Range("E38").Value = '  +2.45%  '

$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.56'
$ws.Range("E39").Value = '  +6.20%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +3.83%  '

$ws.Range("E41").Value = '  +2.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("E42").Value = '  +3.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.054.44'
$ws.Range("E43").Value = '  -1.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.46'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("E45").Value = '  +2.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.16'
$ws.Range("E46").Value = '  +2.07%  '

$ws.Range("E47").Value = '  +5.71%  '

$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("E49").Value = '  +3.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.32'
$ws.Range("E50").Value = '  +6.71%  '

$ws.Range("E51").Value = '  +2.22%  '
